$d = $word.ActiveDocument

# Update the course header/title text for the new academic year:
#   "EG-247/EG-3068 Signals and Systems 2021-2022"
# becomes
#   "EG-247 Signals and Systems 2022-2023"
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found = $find.Execute("EG-247/EG-3068 Signals and Systems 2021-2022", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "EG-247 Signals and Systems 2022-2023", 2)

if (-not $found) {
    Write-Output "Warning: title text to replace was not found"
}
